$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The numeric-looking stats in C2:F5 are stored as text in this workbook
# (see ignoredErrors/numberStoredAsText in the original file), so force
# the target range to Text format before writing the new values. This
# stops Excel from auto-converting the assigned strings into numbers.
$rng = $ws.Range("C2:F5")
$rng.NumberFormat = "@"

$ws.Range("C2").Value = "39"
$ws.Range("D2").Value = "32"
$ws.Range("E2").Value = "4"
$ws.Range("F2").Value = "1"

$ws.Range("C3").Value = "87"
$ws.Range("D3").Value = "45"
$ws.Range("E3").Value = "12"
$ws.Range("F3").Value = "2"

$ws.Range("C4").Value = "58"
$ws.Range("D4").Value = "45"
$ws.Range("E4").Value = "7"
$ws.Range("F4").Value = "1"

$ws.Range("C5").Value = "30"
$ws.Range("D5").Value = "31"
$ws.Range("E5").Value = "1"
$ws.Range("F5").Value = "1"
